$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting existing rows 240:279 down to 241:280
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new record's data
$ws.Cells.Item(240, 1).Value = 4
$ws.Cells.Item(240, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(240, 3).Value = "Los Lagos"
$ws.Cells.Item(240, 4).Value = 44504
$ws.Cells.Item(240, 5).Value = 10
$ws.Cells.Item(240, 6).Value = 100112006
$ws.Cells.Item(240, 7).Value = "Repollo"
$ws.Cells.Item(240, 8).Value = "Crespo record"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 500
$ws.Cells.Item(240, 11).Value = 1200
$ws.Cells.Item(240, 12).Value = 1200
$ws.Cells.Item(240, 13).Value = 1200
$ws.Cells.Item(240, 14).Value = "`$/unidad"
$ws.Cells.Item(240, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(240, 16).Value = 1200
$ws.Cells.Item(240, 17).Value = 1
$ws.Cells.Item(240, 18).Value = "Hortaliza"
